# UC003 - Cancelar Solicitação de Diária -- "1.2.4 to 1.2.5" + minor text fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Version: 1.0 -> 1.2.5  (D2)
$ws.Range("D2").Value = "1.2.5"

# 2) Precondition text: fix "usuario" -> "usuário" and add trailing period.
#    Repeated in every test-case block (B8, B16, B25, B32, B40).
$precondition = "O usuário devidamente autenticado e na tela inicial de cancelar diárias."
$ws.Range("B8").Value = $precondition
$ws.Range("B16").Value = $precondition
$ws.Range("B25").Value = $precondition
$ws.Range("B32").Value = $precondition
$ws.Range("B40").Value = $precondition

# 3) MSG102 confirmation text: add trailing period.
#    Repeated at D10, D18, D34, D42.
$msg102 = "SYSTEM Exibe a mensagem (MSG102 - Confirmar cancelamento)."
$ws.Range("D10").Value = $msg102
$ws.Range("D18").Value = $msg102
$ws.Range("D34").Value = $msg102
$ws.Range("D42").Value = $msg102

# 4) MSG217 error text: remove stray tab character before the closing paren.
#    Repeated at D11, D35.
$msg217 = "SYSTEM Identifica que o usuário não informou uma justificativa para o cancelamento. Não efetiva o cancelamento e exibe mensagem de erro (MSG217 - Necessário informar uma justificativa para o cancelamento de solicitações) para o usuário."
$ws.Range("D11").Value = $msg217
$ws.Range("D35").Value = $msg217

# 5) MSG205 error text: fix typo "Solcitação" -> "Solicitação".
$msg205 = "SYSTEM Identifica que a solicitação de diária está em situação diferente de 'SOLICITADA PARA EMPENHO' ou 'SOLICITADA PARA PRESTAÇÃO DE CONTAS'.  Impede o cancelamento e exibe mensagem de erro (MSG205 - Solicitação de diária não pode ser cancelada) para o usuário."
$ws.Range("D27").Value = $msg205
